$wb = $excel.ActiveWorkbook

$wsMeta  = $wb.Worksheets.Item("Metadata")
$wsCore  = $wb.Worksheets.Item("Core")
$wsChart = $wb.Worksheets.Item("Test Chart")

# --- Sheet "Test Chart": the code for the "DateTime" chart question used to be a
#     placeholder ("testchartcode0"); give it its real code, matching how the other
#     "Core" question codes were fixed up below. Before changing Metadata!G8's own
#     formatting, copy that distinct format onto Test Chart!A2, mirroring the diff. ---
$wsMeta.Range("G8").Copy()
$wsChart.Range("A2").PasteSpecial(-4122)
$wsChart.Range("A2").Value = "PatientChartingDate"

# --- Sheet "Metadata": normalize G8's formatting back to the common style used by
#     every other cell on the sheet (the same style as G7). ---
$wsMeta.Range("G7").Copy()
$wsMeta.Range("G8").PasteSpecial(-4122)

# --- Sheet "Core": replace the placeholder "code" values with the real question codes
#     (each one equal to the question's "type" value in column B). ---
$wsCore.Range("A2").Value = "ComplexChartInstanceName"
$wsCore.Range("A3").Value = "ComplexChartDate"
$wsCore.Range("A4").Value = "ComplexChartType"
$wsCore.Range("A5").Value = "ComplexChartSubtype"
